$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that starts with "- Caloric materials will have"
# ------------------------------------------------------------------
$rParaStart = $d.Content
$rParaStart.Find.ClearFormatting()
$okStart = $rParaStart.Find.Execute("- Caloric materials will have", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okStart) {
    throw "Could not find the target paragraph start"
}
$paraStart = $rParaStart.Start

# Locate the second "field" occurrence - the one immediately followed by
# ". If necessary" (this is the run that also carries lastRenderedPageBreak).
$rField = $d.Content
$rField.Find.ClearFormatting()
$okField = $rField.Find.Execute("field. If necessary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okField) {
    throw "Could not find 'field. If necessary' anchor"
}
$fieldStart = $rField.Start
$fieldEnd = $fieldStart + 5   # length of the word "field"

# ------------------------------------------------------------------
# Step 1: merge "- Caloric materials will have" / "adiabatic temperature
# change" / ", which must be zero ... external " into a single run, by
# doing a no-op Find/Replace restricted to that sub-range only (this
# leaves the following run, which starts at $fieldStart, untouched).
# ------------------------------------------------------------------
$rMerge = $d.Range($paraStart, $fieldStart)
$rMerge.Find.ClearFormatting()
$rMerge.Find.Replacement.ClearFormatting()
$rMerge.Find.Execute("will have adiabatic", $true, $false, $false, $false, $false, $true, 1, $false, "will have adiabatic", 2)

# ------------------------------------------------------------------
# Step 2: insert the new text ", with the correct sign" right after the
# word "field" (and before the following ". If necessary ...").
# ------------------------------------------------------------------
$rWord = $d.Range($fieldStart, $fieldEnd)
$rWord.InsertAfter(", with the correct sign")

# ------------------------------------------------------------------
# Step 3: force "field" and the newly inserted ", with the correct sign"
# to live in their own runs (rather than being silently re-merged with
# neighboring text) by toggling a character-formatting property on then
# back off. This mirrors the run split seen in the target revision.
# ------------------------------------------------------------------
$signLen = ", with the correct sign".Length

$rFieldOnly = $d.Range($fieldStart, $fieldEnd)
$rFieldOnly.Font.Bold = $true
$rFieldOnly.Font.Bold = $false

$rSignOnly = $d.Range($fieldEnd, $fieldEnd + $signLen)
$rSignOnly.Font.Bold = $true
$rSignOnly.Font.Bold = $false
